# Auto-generated edit script: updates crypto price/volume table cells
# to match the target snapshot (commit: "Updated cryptos list on
# Sun Dec  3 02:50:52 UTC 2023 with GitHub Actions").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Forces the cell to hold the exact literal string (so values like
    # "16.00" or "0.400" keep their trailing zeros instead of being
    # auto-converted to a number by Excel), then restores the default
    # "Normal" style so no stray number-format/style is left behind.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '39.426.72'
$ws.Range('E2').Value = '  +1.86%  '
$ws.Range('D3').Value = '2.170.97'
$ws.Range('E3').Value = '  +3.84%  '
$ws.Range('E4').Value = '  +0.04%  '
Set-TextValue 'D5' '230.16'
$ws.Range('E5').Value = '  +0.50%  '
$ws.Range('E6').Value = '  +1.28%  '
Set-TextValue 'D7' '65.16'
$ws.Range('E7').Value = '  +6.76%  '
$ws.Range('E8').Value = '  +0.07%  '
Set-TextValue 'D9' '0.400'
$ws.Range('E9').Value = '  +3.98%  '
Set-TextValue 'D10' '0.0864'
$ws.Range('E10').Value = '  +2.47%  '
$ws.Range('E11').Value = '  +0.01%  '
Set-TextValue 'D12' '16.00'
$ws.Range('E12').Value = '  +5.40%  '
$ws.Range('D13').Value = '2.492.88'
$ws.Range('E13').Value = '  +3.80%  '
Set-TextValue 'D14' '22.51'
$ws.Range('E14').Value = '  +2.33%  '
$ws.Range('E15').Value = '  +0.32%  '
Set-TextValue 'D16' '5.59'
$ws.Range('E16').Value = '  +2.19%  '
$ws.Range('D17').Value = '2.161.64'
$ws.Range('E17').Value = '  +3.70%  '
$ws.Range('D18').Value = '39.446.51'
$ws.Range('E18').Value = '  +2.05%  '
$ws.Range('E19').Value = '  +1.74%  '
Set-TextValue 'D20' '72.40'
$ws.Range('E20').Value = '  +1.01%  '
$ws.Range('D21').Value = '0.0₃0856'
$ws.Range('E21').Value = '  +1.94%  '
Set-TextValue 'D22' '232.52'
$ws.Range('E22').Value = '  +2.30%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('E24').Value = '  -1.22%  '
$ws.Range('E25').Value = '  +2.06%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D26' '172.69'
$ws.Range('E26').Value = '  +1.03%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D27' '9.59'
$ws.Range('E27').Value = '  +0.43%  '
$ws.Range('E28').Value = '  -0.44%  '
Set-TextValue 'D29' '20.09'
$ws.Range('E29').Value = '  +4.44%  '
$ws.Range('E30').Value = '  -1.58%  '
$ws.Range('E31').Value = '  +12.17%  '
Set-TextValue 'D32' '0.122'
$ws.Range('E32').Value = '  +1.90%  '
Set-TextValue 'D33' '4.66'
$ws.Range('E33').Value = '  +3.42%  '
$ws.Range('E34').Value = '  +3.08%  '
Set-TextValue 'D35' '7.17'
$ws.Range('E35').Value = '  +9.47%  '
Set-TextValue 'D36' '0.0621'
$ws.Range('E36').Value = '  +1.84%  '
Set-TextValue 'D37' '2.45'
$ws.Range('E37').Value = '  +2.17%  '
$ws.Range('E38').Value = '  +0.84%  '
$ws.Range('E39').Value = '  -0.11%  '
Set-TextValue 'D40' '104.79'
$ws.Range('E40').Value = '  +3.92%  '
Set-TextValue 'D42' '17.95'
$ws.Range('E42').Value = '  -0.13%  '
$ws.Range('D43').Value = '1.542.16'
$ws.Range('E43').Value = '  +0.54%  '
$ws.Range('E44').Value = '  +5.52%  '
Set-TextValue 'D45' '4.37'
$ws.Range('E45').Value = '  +7.15%  '
Set-TextValue 'D46' '7.99'
$ws.Range('E46').Value = '  +2.97%  '
$ws.Range('E47').Value = '  +7.89%  '
$ws.Range('E48').Value = '  +1.26%  '
$ws.Range('E49').Value = '  +0.46%  '
$ws.Range('D50').Value = '2.374.91'
$ws.Range('E50').Value = '  +3.80%  '
$ws.Range('E51').Value = '  +0.35%  '
